# Update the "Volume(1h)" percentage column (E2:E51) with refreshed values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2  = "  +1.28%  "
    3  = "  +1.64%  "
    4  = "  +0.29%  "
    5  = "  +1.36%  "
    6  = "  +0.24%  "
    7  = "  +1.34%  "
    8  = "  +2.73%  "
    9  = "  +1.08%  "
    10 = "  +2.24%  "
    11 = "  +5.94%  "
    12 = "  +1.19%  "
    13 = "  -1.82%  "
    14 = "  +0.95%  "
    15 = "  +1.87%  "
    16 = "  +0.52%  "
    17 = "  +0.26%  "
    18 = "  +1.00%  "
    19 = "  +0.23%  "
    20 = "  +1.23%  "
    21 = "  +2.88%  "
    22 = "  +0.78%  "
    23 = "  +0.90%  "
    24 = "  +0.17%  "
    25 = "  -0.08%  "
    26 = "  +1.85%  "
    27 = "  +0.04%  "
    28 = "  +1.37%  "
    29 = "  +1.21%  "
    30 = "  +1.18%  "
    31 = "  +5.37%  "
    32 = "  +3.65%  "
    33 = "  +4.53%  "
    34 = "  +1.70%  "
    35 = "  +1.63%  "
    36 = "  +1.39%  "
    37 = "  +3.93%  "
    38 = "  +1.09%  "
    39 = "  +2.24%  "
    40 = "  +2.86%  "
    41 = "  +1.06%  "
    42 = "  +1.27%  "
    43 = "  +2.51%  "
    44 = "  +4.79%  "
    45 = "  +0.32%  "
    46 = "  +0.29%  "
    47 = "  +2.92%  "
    48 = "  +2.13%  "
    49 = "  +2.93%  "
    50 = "  +0.28%  "
    51 = "  +4.44%  "
}

foreach ($row in $values.Keys) {
    $ws.Range("E$row").Value = $values[$row]
}
